$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (Changed) date column C for every existing data
#    row (2..428) from 45192 to 45202.
$ws.Range("C2:C428").Value = 45202

# 2) Row 428 gains an explicit row height (it previously had none, now
#    matches the other data rows at 15pt).
$ws.Rows.Item(428).RowHeight = 15

# 3) Append new row 429: "A 47036-2023"
$ws.Cells.Item(429, 1).Value = "A 47036-2023"
$ws.Cells.Item(429, 2).Value = 45196
$ws.Cells.Item(429, 3).Value = 45202
$ws.Cells.Item(429, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(429, 5).Value = "MORA"
$ws.Cells.Item(429, 7).Value = 2.7
$ws.Cells.Item(429, 8).Value = 0
$ws.Cells.Item(429, 9).Value = 0
$ws.Cells.Item(429, 10).Value = 0
$ws.Cells.Item(429, 11).Value = 0
$ws.Cells.Item(429, 12).Value = 0
$ws.Cells.Item(429, 13).Value = 0
$ws.Cells.Item(429, 14).Value = 0
$ws.Cells.Item(429, 15).Value = 0
$ws.Cells.Item(429, 16).Value = 0
$ws.Cells.Item(429, 17).Value = 0
$ws.Cells.Item(429, 18).Value = ""
$ws.Cells.Item(429, 18).WrapText = $true
$ws.Rows.Item(429).RowHeight = 15

# 4) Append new row 430: "A 46612-2023"
$ws.Cells.Item(430, 1).Value = "A 46612-2023"
$ws.Cells.Item(430, 2).Value = 45198
$ws.Cells.Item(430, 3).Value = 45202
$ws.Cells.Item(430, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(430, 5).Value = "MORA"
$ws.Cells.Item(430, 7).Value = 3.2
$ws.Cells.Item(430, 8).Value = 0
$ws.Cells.Item(430, 9).Value = 0
$ws.Cells.Item(430, 10).Value = 0
$ws.Cells.Item(430, 11).Value = 0
$ws.Cells.Item(430, 12).Value = 0
$ws.Cells.Item(430, 13).Value = 0
$ws.Cells.Item(430, 14).Value = 0
$ws.Cells.Item(430, 15).Value = 0
$ws.Cells.Item(430, 16).Value = 0
$ws.Cells.Item(430, 17).Value = 0
$ws.Cells.Item(430, 18).Value = ""
$ws.Cells.Item(430, 18).WrapText = $true

# Date columns B and C should carry the same date formatting as the rest
# of the table (style copied from an existing formatted date cell).
$dateFormat = $ws.Range("B428").NumberFormat
$ws.Range("B429:C430").NumberFormat = $dateFormat
